# "Add Filter tahun dan bulan"
# Adds a new column M ("telkomsel") to both data rows and changes C2's
# value from the number 1 to the text "tes" (new filter values used for
# the year/month filter columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header): new cell M1 = "telkomsel"
$ws.Range("M1").Value = "telkomsel"

# Row 2 (data): C2 changes from numeric 1 to text "tes"
$ws.Range("C2").Value = "tes"

# Row 2 (data): new cell M2 = "telkomsel"
$ws.Range("M2").Value = "telkomsel"

# Move/collapse the selection onto C2, matching the saved view state
[void]$ws.Range("C2").Select()
